$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.503.62"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.924.39"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.01"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  -1.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2878"
$ws.Range("E8").Value = "  -2.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06748"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "106.37"
$ws.Range("E10").Value = "  +5.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.26"
$ws.Range("E11").Value = "  -2.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07756"
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.903.16"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.295"
$ws.Range("E14").Value = "  +2.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6602"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "291.45"
$ws.Range("E16").Value = "  -5.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.495.84"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007593"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.93"
$ws.Range("E20").Value = "  -2.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.152.85"
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.269"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.206"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.366"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.74"
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.38"
$ws.Range("E27").Value = "  +4.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.104"
$ws.Range("E28").Value = "  +7.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1070"
$ws.Range("E29").Value = "  -5.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.365"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.178"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.996"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05039"
$ws.Range("E33").Value = "  -1.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7431"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02117"
$ws.Range("E36").Value = "  +7.55%  "
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.683"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.082"
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.21"
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8742"
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.858"
$ws.Range("E42").Value = "  +3.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4274"
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "67.43"
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "49.51"
$ws.Range("E46").Value = "  +15.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.196"
$ws.Range("E47").Value = "  -2.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.264"
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1218"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2458"
$ws.Range("E51").Value = "  +9.29%  "
